# Fruta / hortaliza, semanal
# Update Fecha (D) and Volumen (M) values across rows 2-9, and swap the
# Unidad de comercialización / Precio $/Kg / Kg por unidad (Q/S/T) values
# between rows 3 and 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) updates
$ws.Range("D2").Value = 44330
$ws.Range("D3").Value = 44306
$ws.Range("D4").Value = 44302
$ws.Range("D5").Value = 44323
$ws.Range("D6").Value = 44309
$ws.Range("D7").Value = 44327
$ws.Range("D8").Value = 44313
$ws.Range("D9").Value = 44322

# Column M (Volumen) updates
$ws.Range("M2").Value = 60
$ws.Range("M4").Value = 80
$ws.Range("M6").Value = 80
$ws.Range("M7").Value = 60
$ws.Range("M8").Value = 120

# Rows 3 and 6: swap Q (Unidad de comercialización), S (Precio $/Kg), T (Kg / unidad)
$ws.Range("Q3").Value = "$/caja 10 kilos empedrada"
$ws.Range("S3").Value = 11500
$ws.Range("T3").Value = 1

$ws.Range("Q6").Value = "$/caja 14 kilos granel"
$ws.Range("S6").Value = 821
$ws.Range("T6").Value = 14
